$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update time_taken (column F) timestamps on the "data" sheet ---
$ws1.Range("F2").Value = "2021-10-05 14:34:32.196960"
$ws1.Range("F3").Value = "2021-10-05 14:34:32.196968"
$ws1.Range("F4").Value = "2021-10-05 14:34:32.196972"
$ws1.Range("F5").Value = "2021-10-05 14:34:32.196974"
$ws1.Range("F6").Value = "2021-10-05 14:34:32.196977"
$ws1.Range("F7").Value = "2021-10-05 14:34:32.196980"
$ws1.Range("F8").Value = "2021-10-05 14:34:32.196983"
$ws1.Range("F9").Value = "2021-10-05 14:34:32.196985"
$ws1.Range("F10").Value = "2021-10-05 14:34:32.196988"
$ws1.Range("F11").Value = "2021-10-05 14:34:32.196991"
$ws1.Range("F12").Value = "2021-10-05 14:34:32.196993"
$ws1.Range("F13").Value = "2021-10-05 14:34:32.196996"
$ws1.Range("F14").Value = "2021-10-05 14:34:32.196998"
$ws1.Range("F15").Value = "2021-10-05 14:34:32.197001"
$ws1.Range("F16").Value = "2021-10-05 14:34:32.197003"
$ws1.Range("F17").Value = "2021-10-05 14:34:32.197006"
$ws1.Range("F18").Value = "2021-10-05 14:34:32.197009"
$ws1.Range("F19").Value = "2021-10-05 14:34:32.197012"
$ws1.Range("F20").Value = "2021-10-05 14:34:32.197014"
$ws1.Range("F21").Value = "2021-10-05 14:34:32.197017"
$ws1.Range("F22").Value = "2021-10-05 14:34:32.197019"
$ws1.Range("F23").Value = "2021-10-05 14:34:32.197022"
$ws1.Range("F24").Value = "2021-10-05 14:34:32.197024"
$ws1.Range("F25").Value = "2021-10-05 14:34:32.197027"
$ws1.Range("F26").Value = "2021-10-05 14:34:32.197030"
$ws1.Range("F27").Value = "2021-10-05 14:34:32.197032"
$ws1.Range("F28").Value = "2021-10-05 14:34:32.197035"
$ws1.Range("F29").Value = "2021-10-05 14:34:32.197037"
$ws1.Range("F30").Value = "2021-10-05 14:34:32.197040"
$ws1.Range("F31").Value = "2021-10-05 14:34:32.197042"
$ws1.Range("F32").Value = "2021-10-05 14:34:32.197045"
$ws1.Range("F33").Value = "2021-10-05 14:34:32.197047"
$ws1.Range("F34").Value = "2021-10-05 14:34:32.197050"
$ws1.Range("F35").Value = "2021-10-05 14:34:32.197053"
$ws1.Range("F36").Value = "2021-10-05 14:34:32.197056"
$ws1.Range("F37").Value = "2021-10-05 14:34:32.197058"
$ws1.Range("F38").Value = "2021-10-05 14:34:32.197061"
$ws1.Range("F39").Value = "2021-10-05 14:34:32.197063"
$ws1.Range("F40").Value = "2021-10-05 14:34:32.197066"
$ws1.Range("F41").Value = "2021-10-05 14:34:32.197068"
$ws1.Range("F42").Value = "2021-10-05 14:34:32.197071"
$ws1.Range("F43").Value = "2021-10-05 14:34:32.197074"
$ws1.Range("F44").Value = "2021-10-05 14:34:32.197076"
$ws1.Range("F45").Value = "2021-10-05 14:34:32.197079"
$ws1.Range("F46").Value = "2021-10-05 14:34:32.197081"
$ws1.Range("F47").Value = "2021-10-05 14:34:32.197084"
$ws1.Range("F48").Value = "2021-10-05 14:34:32.197086"
$ws1.Range("F49").Value = "2021-10-05 14:34:32.197089"
$ws1.Range("F50").Value = "2021-10-05 14:34:32.197091"
$ws1.Range("F51").Value = "2021-10-05 14:34:32.197094"
$ws1.Range("F52").Value = "2021-10-05 14:34:32.197096"
$ws1.Range("F53").Value = "2021-10-05 14:34:32.197099"
$ws1.Range("F54").Value = "2021-10-05 14:34:32.197102"
$ws1.Range("F55").Value = "2021-10-05 14:34:32.197104"
$ws1.Range("F56").Value = "2021-10-05 14:34:32.197107"
$ws1.Range("F57").Value = "2021-10-05 14:34:32.197109"
$ws1.Range("F58").Value = "2021-10-05 14:34:32.197112"
$ws1.Range("F59").Value = "2021-10-05 14:34:32.197114"
$ws1.Range("F60").Value = "2021-10-05 14:34:32.197117"
$ws1.Range("F61").Value = "2021-10-05 14:34:32.197120"
$ws1.Range("F62").Value = "2021-10-05 14:34:32.197122"
$ws1.Range("F63").Value = "2021-10-05 14:34:32.197125"
$ws1.Range("F64").Value = "2021-10-05 14:34:32.197127"
$ws1.Range("F65").Value = "2021-10-05 14:34:32.197130"
$ws1.Range("F66").Value = "2021-10-05 14:34:32.197134"
$ws1.Range("F67").Value = "2021-10-05 14:34:32.197137"
$ws1.Range("F68").Value = "2021-10-05 14:34:32.197139"
$ws1.Range("F69").Value = "2021-10-05 14:34:32.197142"
$ws1.Range("F70").Value = "2021-10-05 14:34:32.197144"
$ws1.Range("F71").Value = "2021-10-05 14:34:32.197147"
$ws1.Range("F72").Value = "2021-10-05 14:34:32.197149"
$ws1.Range("F73").Value = "2021-10-05 14:34:32.197152"
$ws1.Range("F74").Value = "2021-10-05 14:34:32.197155"
$ws1.Range("F75").Value = "2021-10-05 14:34:32.197157"
$ws1.Range("F76").Value = "2021-10-05 14:34:32.197159"
$ws1.Range("F77").Value = "2021-10-05 14:34:32.197162"
$ws1.Range("F78").Value = "2021-10-05 14:34:32.197167"
$ws1.Range("F79").Value = "2021-10-05 14:34:32.197170"
$ws1.Range("F80").Value = "2021-10-05 14:34:32.197172"
$ws1.Range("F81").Value = "2021-10-05 14:34:32.197175"
$ws1.Range("F82").Value = "2021-10-05 14:34:32.197178"
$ws1.Range("F83").Value = "2021-10-05 14:34:32.197180"
$ws1.Range("F84").Value = "2021-10-05 14:34:32.197183"
$ws1.Range("F85").Value = "2021-10-05 14:34:32.197185"
$ws1.Range("F86").Value = "2021-10-05 14:34:32.197187"
$ws1.Range("F87").Value = "2021-10-05 14:34:32.197190"
$ws1.Range("F88").Value = "2021-10-05 14:34:32.197193"
$ws1.Range("F89").Value = "2021-10-05 14:34:32.197195"
$ws1.Range("F90").Value = "2021-10-05 14:34:32.197198"
$ws1.Range("F91").Value = "2021-10-05 14:34:32.197200"
$ws1.Range("F92").Value = "2021-10-05 14:34:32.197203"
$ws1.Range("F93").Value = "2021-10-05 14:34:32.197205"
$ws1.Range("F94").Value = "2021-10-05 14:34:32.197209"
$ws1.Range("F95").Value = "2021-10-05 14:34:32.197212"
$ws1.Range("F96").Value = "2021-10-05 14:34:32.197214"
$ws1.Range("F97").Value = "2021-10-05 14:34:32.197217"
$ws1.Range("F98").Value = "2021-10-05 14:34:32.197219"
$ws1.Range("F99").Value = "2021-10-05 14:34:32.197222"
$ws1.Range("F100").Value = "2021-10-05 14:34:32.197224"
$ws1.Range("F101").Value = "2021-10-05 14:34:32.197227"
$ws1.Range("F102").Value = "2021-10-05 14:34:32.197229"
$ws1.Range("F103").Value = "2021-10-05 14:34:32.197232"
$ws1.Range("F104").Value = "2021-10-05 14:34:32.197234"
$ws1.Range("F105").Value = "2021-10-05 14:34:32.197237"
$ws1.Range("F106").Value = "2021-10-05 14:34:32.197239"
$ws1.Range("F107").Value = "2021-10-05 14:34:32.197242"
$ws1.Range("F108").Value = "2021-10-05 14:34:32.197244"
$ws1.Range("F109").Value = "2021-10-05 14:34:32.197247"
$ws1.Range("F110").Value = "2021-10-05 14:34:32.197252"
$ws1.Range("F111").Value = "2021-10-05 14:34:32.197255"
$ws1.Range("F112").Value = "2021-10-05 14:34:32.197258"
$ws1.Range("F113").Value = "2021-10-05 14:34:32.197260"
$ws1.Range("F114").Value = "2021-10-05 14:34:32.197263"
$ws1.Range("F115").Value = "2021-10-05 14:34:32.197265"
$ws1.Range("F116").Value = "2021-10-05 14:34:32.197268"
$ws1.Range("F117").Value = "2021-10-05 14:34:32.197270"
$ws1.Range("F118").Value = "2021-10-05 14:34:32.197273"
$ws1.Range("F119").Value = "2021-10-05 14:34:32.197275"
$ws1.Range("F120").Value = "2021-10-05 14:34:32.197278"
$ws1.Range("F121").Value = "2021-10-05 14:34:32.197280"
$ws1.Range("F122").Value = "2021-10-05 14:34:32.197283"
$ws1.Range("F123").Value = "2021-10-05 14:34:32.197285"
$ws1.Range("F124").Value = "2021-10-05 14:34:32.197288"
$ws1.Range("F125").Value = "2021-10-05 14:34:32.197290"
$ws1.Range("F126").Value = "2021-10-05 14:34:32.197293"
$ws1.Range("F127").Value = "2021-10-05 14:34:32.197296"
$ws1.Range("F128").Value = "2021-10-05 14:34:32.197298"
$ws1.Range("F129").Value = "2021-10-05 14:34:32.197300"
$ws1.Range("F130").Value = "2021-10-05 14:34:32.197305"
$ws1.Range("F131").Value = "2021-10-05 14:34:32.197308"
$ws1.Range("F132").Value = "2021-10-05 14:34:32.197311"
$ws1.Range("F133").Value = "2021-10-05 14:34:32.197313"
$ws1.Range("F134").Value = "2021-10-05 14:34:32.197316"

# --- Add new "metadata" sheet, placed after "data" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "metadata"

# Header row (bold, bordered, centered, top-aligned -- matches "data" header style)
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

$headerRng = $ws2.Range("B1:G1")
$headerRng.Font.Bold = $true
$headerRng.Borders.LineStyle = 1
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160

# Data row 2
$ws2.Range("A2").Value = 0
$ws2.Range("A2").Font.Bold = $true
$ws2.Range("A2").Borders.LineStyle = 1
$ws2.Range("A2").HorizontalAlignment = -4108
$ws2.Range("A2").VerticalAlignment = -4160

$ws2.Range("B2").Value = "Macrocephaly_Megalencephaly"
$ws2.Range("C2").Value = 135

# Force D2 to be stored as TEXT ("0.88"), not a number, then drop the
# number-format override so the cell keeps the default (unstyled) look.
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "0.88"
$ws2.Range("D2").ClearFormats()

$ws2.Range("E2").Value = "2021-10-04T00:57:54.786887Z"
$ws2.Range("F2").Value = "2021-10-05 14:34:32.193628"
$ws2.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/135/?format=json"

$wb.Worksheets.Item(1).Select()
